# Update "Förändrad" (column C) date for every data row (2..176): 45189 -> 45190
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 176; $r++) {
    $ws.Cells.Item($r, 3).Value = 45190
}

# Row 2 (A 14918-2023): signal species counts grew (3 new species found)
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 4
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 40

$crlf = "`r`n"
$row2Species = @(
    "Knärot",
    "Kopparspindling",
    "Lakritsmusseron",
    "Läderdoftande fingersvamp",
    "Dofttaggsvamp",
    "Flattoppad klubbsvamp",
    "Gultoppig fingersvamp",
    "Odörspindling",
    "Orange taggsvamp",
    "Persiljespindling",
    "Svartvit taggsvamp",
    "Talltita",
    "Tretåig hackspett",
    "Ullticka",
    "Vedtrappmossa",
    "Äggvaxskivling",
    "Anisspindling",
    "Bronshjon",
    "Fjällig taggsvamp s.str.",
    "Korallblylav",
    "Mörk husmossa",
    "Olivspindling",
    "Rödgul trumpetsvamp",
    "Skarp dropptaggsvamp",
    "Skinnlav",
    "Skogshakmossa",
    "Sotriska",
    "Spindelblomster",
    "Svart trolldruva",
    "Svavelriska",
    "Thomsons trägnagare",
    "Tibast",
    "Tvåblad",
    "Underviol",
    "Vågbandad barkbock",
    "Vårärt",
    "Zontaggsvamp",
    "Fläcknycklar",
    "Blåsippa",
    "Lopplummer"
)
$ws.Range("R2").Value = $row2Species -join $crlf

# Row 10 (A 13750-2023): one new species found (Leptoporus mollis)
$ws.Range("J10").Value = 3
$ws.Range("O10").Value = 3
$ws.Range("Q10").Value = 5

$row10Species = @(
    "Gultoppig fingersvamp",
    "Leptoporus mollis",
    "Talltita",
    "Fjällig taggsvamp s.str.",
    "Vedticka"
)
$ws.Range("R10").Value = $row10Species -join $crlf
